$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update cell values. Cells are pre-formatted as Text before the write
# (then format is cleared back to default/General) so that numeric-looking
# strings like "1.003" or "24.924.80" are preserved verbatim as text,
# matching the source data which stores every Price/Volume cell as a string.
function Set-TextValue($rangeAddr, $val) {
    $rng = $ws.Range($rangeAddr)
    $rng.NumberFormat = "@"
    $rng.Value = $val
    $rng.ClearFormats()
}

Set-TextValue 'D2' '24.924.80'
Set-TextValue 'E2' '  +2.18%  '
Set-TextValue 'D3' '1.709.17'
Set-TextValue 'E3' '  +1.88%  '
Set-TextValue 'D4' '1.003'
Set-TextValue 'E4' '  -0.36%  '
Set-TextValue 'D5' '313.81'
Set-TextValue 'E5' '  +2.25%  '
Set-TextValue 'D6' '0.9985'
Set-TextValue 'E6' '  -0.16%  '
Set-TextValue 'D7' '0.3750'
Set-TextValue 'E7' '  +1.22%  '
Set-TextValue 'D8' '49.41'
Set-TextValue 'E8' '  +3.60%  '
Set-TextValue 'D9' '0.3448'
Set-TextValue 'E9' '  +0.43%  '
Set-TextValue 'D10' '1.226'
Set-TextValue 'E10' '  +5.12%  '
Set-TextValue 'D11' '0.07563'
Set-TextValue 'E11' '  +4.36%  '
Set-TextValue 'D12' '0.9992'
Set-TextValue 'E12' '  -0.40%  '
Set-TextValue 'D13' '21.31'
Set-TextValue 'E13' '  +5.75%  '
Set-TextValue 'D14' '6.331'
Set-TextValue 'E14' '  +3.93%  '
Set-TextValue 'D15' '7.092'
Set-TextValue 'E15' '  +5.38%  '
Set-TextValue 'D16' '1.709.03'
Set-TextValue 'E16' '  +1.73%  '
Set-TextValue 'D17' '0.00001133'
Set-TextValue 'E17' '  +2.60%  '
Set-TextValue 'D18' '0.06731'
Set-TextValue 'E18' '  +0.98%  '
Set-TextValue 'D19' '0.9978'
Set-TextValue 'E19' '  -0.29%  '
Set-TextValue 'D20' '84.15'
Set-TextValue 'E20' '  +3.91%  '
Set-TextValue 'D21' '17.35'
Set-TextValue 'E21' '  +5.79%  '
Set-TextValue 'D22' '6.407'
Set-TextValue 'E22' '  +5.07%  '
Set-TextValue 'D23' '13.11'
Set-TextValue 'E23' '  +8.08%  '
Set-TextValue 'D24' '24.916.22'
Set-TextValue 'E24' '  +2.29%  '
Set-TextValue 'D25' '2.446'
Set-TextValue 'E25' '  -0.61%  '
Set-TextValue 'D26' '2.804'
Set-TextValue 'E26' '  +6.01%  '
Set-TextValue 'D27' '20.47'
Set-TextValue 'E27' '  +5.35%  '
Set-TextValue 'D28' '149.72'
Set-TextValue 'E28' '  -2.25%  '
Set-TextValue 'D29' '132.98'
Set-TextValue 'E29' '  +4.58%  '
Set-TextValue 'B30' 'ImmutableX'
Set-TextValue 'C30' 'https://coinranking.com/coin/Z96jIvLU7+immutablex-imx'
Set-TextValue 'D30' '1.254'
Set-TextValue 'E30' '  +29.50%  '
Set-TextValue 'B31' 'WrappedliquidstakedEther2.0'
Set-TextValue 'C31' 'https://coinranking.com/coin/CiixT63n3+wrappedliquidstakedether20-wsteth'
Set-TextValue 'D31' '1.896.49'
Set-TextValue 'E31' '  +1.74%  '
Set-TextValue 'D32' '6.841'
Set-TextValue 'E32' '  +9.04%  '
Set-TextValue 'D33' '4.219'
Set-TextValue 'E33' '  +3.96%  '
Set-TextValue 'E34' '  +13.76%  '
Set-TextValue 'B35' 'Stellar'
Set-TextValue 'C35' 'https://coinranking.com/coin/f3iaFeCKEmkaZ+stellar-xlm'
Set-TextValue 'D35' '0.08805'
Set-TextValue 'E35' '  +4.24%  '
Set-TextValue 'B36' 'WEMIXTOKEN'
Set-TextValue 'C36' 'https://coinranking.com/coin/08CsQa-Ov+wemixtoken-wemix'
Set-TextValue 'D36' '1.772'
Set-TextValue 'E36' '  +4.55%  '
Set-TextValue 'D37' '5.658'
Set-TextValue 'E37' '  +6.45%  '
Set-TextValue 'D38' '0.06672'
Set-TextValue 'E38' '  +3.56%  '
Set-TextValue 'D39' '9.200'
Set-TextValue 'E39' '  +3.91%  '
Set-TextValue 'D40' '0.02419'
Set-TextValue 'E40' '  +4.49%  '
Set-TextValue 'D41' '0.2243'
Set-TextValue 'E41' '  +7.54%  '
Set-TextValue 'D42' '1.276'
Set-TextValue 'E42' '  +2.43%  '
Set-TextValue 'D43' '0.6481'
Set-TextValue 'E43' '  +5.72%  '
Set-TextValue 'D44' '0.9979'
Set-TextValue 'E44' '  -0.20%  '
Set-TextValue 'D45' '13.87'
Set-TextValue 'E45' '  +5.21%  '
Set-TextValue 'D46' '0.6171'
Set-TextValue 'E46' '  +4.57%  '
Set-TextValue 'D47' '3.843'
Set-TextValue 'E47' '  +1.94%  '
Set-TextValue 'D48' '2.135'
Set-TextValue 'E48' '  +5.89%  '
Set-TextValue 'D49' '129.76'
Set-TextValue 'E49' '  +2.17%  '
Set-TextValue 'D50' '0.07323'
Set-TextValue 'E50' '  +2.21%  '
Set-TextValue 'D51' '80.34'
Set-TextValue 'E51' '  +6.34%  '
